# Generate Report for Archive
#
# 1. Update the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears (Overview!E2:F2, and the
#    "Status" column (C2) on each per-locale sheet).
# 2. Shrink the two "status" columns' width to match the new,
#    shorter text (17.2159881591797 -> 13.4101845877511 in raw OOXML
#    units, which is a ColumnWidth of ~12.5 once the fixed 5/6
#    padding Excel adds is accounted for).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# New narrower width for the "status" columns, expressed as a
# `ColumnWidth` (character units). Excel stores columns internally as
# ColumnWidth + 5/6, so this lands on the desired ~13.41 raw width.
$newColumnWidth = 12.5

# --- Overview sheet: columns E (zh-cn) and F (de-de) ---
if ($wsOverview.Range("E2").Value() -eq $oldStatus) {
    $wsOverview.Range("E2").Value = $newStatus
}
if ($wsOverview.Range("F2").Value() -eq $oldStatus) {
    $wsOverview.Range("F2").Value = $newStatus
}
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# --- Per-locale sheets: column C is "Status" ---
if ($wsZhCn.Range("C2").Value() -eq $oldStatus) {
    $wsZhCn.Range("C2").Value = $newStatus
}
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

if ($wsDeDe.Range("C2").Value() -eq $oldStatus) {
    $wsDeDe.Range("C2").Value = $newStatus
}
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
